$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2024922118380062
$ws.Range("C2").Value = 0.5700934579439252
$ws.Range("J2").Value = 0.009345794392523364
$ws.Range("P2").Value = 0.1495327102803738
$ws.Range("S2").Value = 0.06853582554517133

# Row 3
$ws.Range("B3").Value = 0.005263157894736842
$ws.Range("C3").Value = 0.03157894736842105
$ws.Range("J3").Value = 0.01052631578947368
$ws.Range("P3").Value = 0.7842105263157895
$ws.Range("S3").Value = 0.1684210526315789

# Row 6
$ws.Range("B6").Value = 0.07421875
$ws.Range("D6").Value = 0.0078125
$ws.Range("E6").Value = 0.0078125
$ws.Range("F6").Value = 0.0546875
$ws.Range("J6").Value = 0.265625
$ws.Range("O6").Value = 0.015625
$ws.Range("Q6").Value = 0.16796875
$ws.Range("R6").Value = 0.09375
$ws.Range("S6").Value = 0.3125

# Row 7
$ws.Range("B7").Value = 0.07746478873239436
$ws.Range("D7").Value = 0.01408450704225352
$ws.Range("F7").Value = 0.09859154929577464
$ws.Range("J7").Value = 0.1549295774647887
$ws.Range("O7").Value = 0.02112676056338028
$ws.Range("Q7").Value = 0.1901408450704225
$ws.Range("R7").Value = 0.06338028169014084
$ws.Range("S7").Value = 0.3802816901408451

# Row 8
$ws.Range("B8").Value = 0.09555555555555556
$ws.Range("D8").Value = 0.02444444444444445
$ws.Range("E8").Value = 0.002222222222222222
$ws.Range("F8").Value = 0.06444444444444444
$ws.Range("J8").Value = 0.1266666666666667
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.1911111111111111
$ws.Range("R8").Value = 0.1044444444444445
$ws.Range("S8").Value = 0.3711111111111111

# Row 9
$ws.Range("B9").Value = 0.115702479338843
$ws.Range("D9").Value = 0.01652892561983471
$ws.Range("F9").Value = 0.06198347107438017
$ws.Range("J9").Value = 0.115702479338843
$ws.Range("O9").Value = 0.02479338842975207
$ws.Range("Q9").Value = 0.2066115702479339
$ws.Range("R9").Value = 0.1322314049586777
$ws.Range("S9").Value = 0.3264462809917356

# Row 10
$ws.Range("B10").Value = 0.1184798807749627
$ws.Range("D10").Value = 0.02384500745156483
$ws.Range("E10").Value = 0.0007451564828614009
$ws.Range("F10").Value = 0.07749627421758569
$ws.Range("J10").Value = 0.1073025335320417
$ws.Range("O10").Value = 0.01937406855439642
$ws.Range("Q10").Value = 0.2220566318926975
$ws.Range("R10").Value = 0.08420268256333831
$ws.Range("S10").Value = 0.3464977645305514

# Row 11
$ws.Range("G11").Value = 0.1497797356828194
$ws.Range("J11").Value = 0.1233480176211454
$ws.Range("K11").Value = 0.1982378854625551
$ws.Range("L11").Value = 0.5198237885462555
$ws.Range("S11").Value = 0.00881057268722467

# Row 12
$ws.Range("G12").Value = 0.6929133858267716
$ws.Range("J12").Value = 0.2047244094488189
$ws.Range("K12").Value = 0.007874015748031496
$ws.Range("L12").Value = 0.07086614173228346
$ws.Range("S12").Value = 0.02362204724409449

# Row 15
$ws.Range("F15").Value = 0.02150537634408602
$ws.Range("H15").Value = 0.1577060931899641
$ws.Range("I15").Value = 0.1111111111111111
$ws.Range("J15").Value = 0.3584229390681004
$ws.Range("K15").Value = 0.03225806451612903
$ws.Range("M15").Value = 0.01075268817204301
$ws.Range("N15").Value = 0.003584229390681004
$ws.Range("O15").Value = 0.05734767025089606
$ws.Range("S15").Value = 0.2473118279569892

# Row 16
$ws.Range("F16").Value = 0.03083700440528634
$ws.Range("H16").Value = 0.1541850220264317
$ws.Range("I16").Value = 0.09691629955947137
$ws.Range("J16").Value = 0.4140969162995595
$ws.Range("K16").Value = 0.1013215859030837
$ws.Range("M16").Value = 0.03524229074889868
$ws.Range("O16").Value = 0.06167400881057269
$ws.Range("S16").Value = 0.105726872246696

# Row 17
$ws.Range("F17").Value = 0.01796407185628742
$ws.Range("H17").Value = 0.1996007984031936
$ws.Range("I17").Value = 0.09780439121756487
$ws.Range("J17").Value = 0.4231536926147705
$ws.Range("K17").Value = 0.06786427145708583
$ws.Range("M17").Value = 0.009980039920159681
$ws.Range("N17").Value = 0.001996007984031936
$ws.Range("O17").Value = 0.09181636726546906
$ws.Range("S17").Value = 0.08982035928143713

# Row 18
$ws.Range("F18").Value = 0.004424778761061947
$ws.Range("H18").Value = 0.1858407079646018
$ws.Range("I18").Value = 0.1238938053097345
$ws.Range("J18").Value = 0.495575221238938
$ws.Range("K18").Value = 0.06637168141592921
$ws.Range("M18").Value = 0.008849557522123894
$ws.Range("N18").Value = 0.004424778761061947
$ws.Range("O18").Value = 0.03097345132743363
$ws.Range("S18").Value = 0.07964601769911504

# Row 19
$ws.Range("F19").Value = 0.02642796248934356
$ws.Range("H19").Value = 0.1926683716965047
$ws.Range("I19").Value = 0.09121909633418585
$ws.Range("J19").Value = 0.3938618925831202
$ws.Range("K19").Value = 0.08184143222506395
$ws.Range("M19").Value = 0.01619778346121057
$ws.Range("N19").Value = 0.0008525149190110827
$ws.Range("O19").Value = 0.09633418584825235
$ws.Range("S19").Value = 0.1005967604433078
